# The workbook tracks daily "Pepino ensalada" price records. This edit adds
# a new, more-recent weekly record at the top of the data block (row 899),
# pushing the existing rows 899-1003 down by one (to 900-1004).
#
# Net effect vs. the original file:
#   - sheet dimension grows from A1:R1003 to A1:R1004
#   - a brand-new row of data is inserted at row 899
#   - every previously existing row from 899 through 1003 shifts down by one row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 899 (shifts 899:1003 down to 900:1004,
# inheriting the row-899 formatting for the date column, same as Excel does).
$ws.Rows.Item(899).Insert()

# Populate the newly inserted row 899 with the new record.
$ws.Cells.Item(899, 1).Value  = 8
$ws.Cells.Item(899, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(899, 3).Value  = "Coquimbo"
$ws.Cells.Item(899, 4).Value  = 45212
$ws.Cells.Item(899, 5).Value  = 4
$ws.Cells.Item(899, 6).Value  = 100112043
$ws.Cells.Item(899, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(899, 8).Value  = "Sin especificar"
$ws.Cells.Item(899, 9).Value  = "Primera"
$ws.Cells.Item(899, 10).Value = 560
$ws.Cells.Item(899, 11).Value = 14500
$ws.Cells.Item(899, 12).Value = 15000
$ws.Cells.Item(899, 13).Value = 14750
$ws.Cells.Item(899, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(899, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(899, 16).Value = 246
$ws.Cells.Item(899, 17).Value = 60
$ws.Cells.Item(899, 18).Value = "Hortaliza"
